# Fixed tests with more verbose errors for import; Increased verbosity of
# errors during uploading.
#
# 1) Shared string fix: "Error:not a number;18-03-2022"
#                     -> "Error: not a number;18-03-2022"
# 2) Column widths widened for columns A:J (headers got more verbose, so the
#    best-fit widths grew).
# 3) Cell formatting: columns C,D,E,F,G,I,J on the data rows (2-6) lose their
#    explicit "text" number-format style and fall back to the default style
#    (column H keeps its style).
# 4) View: scrolled right (topLeftCell E1) with I2:J6 selected instead of D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the typo'd error message (missing space after "Error:") -------
$ws.Range("D4").Value = "Error: not a number;18-03-2022"

# --- 2) Resize columns A:J to their new best-fit widths --------------------
# (ColumnWidth is expressed in characters; the stored "width" in the XML is
# ColumnWidth + a fixed padding offset, so we subtract that offset here.)
$colOffset = 0.8333333333333334

$targetWidths = @{
    1  = 9.42578125
    2  = 24
    3  = 18.5703125
    4  = 40
    5  = 28.28515625
    6  = 24.140625
    7  = 23.28515625
    8  = 24.42578125
    9  = 48.85546875
    10 = 73.28515625
}

foreach ($col in 1..10) {
    $ws.Columns($col).ColumnWidth = $targetWidths[$col] - $colOffset
}

# --- 3) Strip the explicit style from the non-header "data" cells so they
#        revert to the default (unstyled) format, matching the default
#        style already used by row 1's header cells.
$defaultStyle = $ws.Range("A1").Style

foreach ($row in 2..6) {
    $ws.Range("C${row}:G${row}").Style = $defaultStyle
    $ws.Range("I${row}:J${row}").Style = $defaultStyle
}

# --- 4) Update the view: scroll so column E is the left-most visible column
#        and select I2:J6 (instead of D5).
$ws.Range("I2:J6").Select()
$excel.ActiveWindow.ScrollColumn = 5
